# "add data for web"
#
# The `dt` column (A) held text timestamps such as "2015-02-30 12:00" /
# "2015-04-60 12:00" — the day-of-month had been incremented past the end
# of the month by a buggy generator instead of always landing on the 15th.
# Fix those rows by entering the *correct* mid-month timestamp as a real
# Excel date/time value (so it becomes a date serial number with a date
# number format), the same way typing a date into a cell in Excel does.
# The rows that already held a valid date string (Jan / Oct, "...-15 12:00")
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> corrected "YYYY-MM-15 12:00" serial date/time value
$fixedDates = @{
  3  = 42050.5   # 2015-02-15 12:00
  4  = 42078.5   # 2015-03-15 12:00
  5  = 42109.5   # 2015-04-15 12:00
  6  = 42139.5   # 2015-05-15 12:00
  7  = 42170.5   # 2015-06-15 12:00
  8  = 42200.5   # 2015-07-15 12:00
  9  = 42231.5   # 2015-08-15 12:00
  10 = 42262.5   # 2015-09-15 12:00
  12 = 42323.5   # 2015-11-15 12:00
  13 = 42353.5   # 2015-12-15 12:00
}

foreach ($row in $fixedDates.Keys) {
  $cell = $ws.Range("A$row")
  $cell.Value = $fixedDates[$row]
  $cell.NumberFormat = "m/d/yy h:mm"
}
